$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove rows 3 and 4 (extra index/json rows no longer present)
$ws.Range("A3:B4").EntireRow.Delete()

# Update the JSON payload in B2 with the new quiz question set
$ws.Range("B2").Value = '{"quiz_questions":[{"id":"1","question":"Name\nAnnie\nBob\nCallie\nDerek\n\nWhich of these query will display the the table given above ?\n","possible_answers":"Select employee from name,Select name,Select name from employee,Select employee","correct_answer":"3"},{"id":"2","question":"Here which of the following displays the unique values of the column?\n\n   SELECT ________ dept_name \n   FROM instructor;\n","possible_answers":"All,From,Distinct,Name","correct_answer":"3"},{"id":"3","question":"This Query can be replaced by which one of the following?\n\n   SELECT name, course_id\n   FROM instructor, teaches\n   WHERE instructor_ID= teaches_ID;\n","possible_answers":"Select name,course_id from teaches,instructor where instructor_id=course_id;,Select name, course_id from instructor natural join teaches;,Select name, course_id from instructor; ,Select course_id from instructor join teaches;","correct_answer":"2"},{"id":"4","question":"SELECT * FROM employee WHERE salary>10000 AND dept_id=101;\nWhich of the following fields are displayed as output?\n","possible_answers":"Salary, dept_id ,Employee ,Salary,All the field of employee relation","correct_answer":"4"},{"id":"5","question":"Which of the following statements contains an error?","possible_answers":"Select * from emp where empid = 10003; ,Select empid from emp where empid = 10006; ,Select empid from emp; ,Select empid where empid = 1009 and lastname = ‘GELLER’;","correct_answer":"4"},{"id":"6","question":"In the given query which of the keyword has to be inserted?\n\nINSERT INTO employee _____ (1002,Joey,2000);\n","possible_answers":"Table ,Values ,Relation ,Field ","correct_answer":"2"} ,{"id":"7","question":"Which of the following makes the transaction permanent in the database? ","possible_answers":"View ,Commit ,Rollback ,Flashback ","correct_answer":"2"} ,{"id":"8","question":"In order to undo the work of transaction after last commit which one should be used? ","possible_answers":"View ,Commit ,Rollback ,Flashback ","correct_answer":"3"} ,{"id":"9","question":"To include integrity constraint in an existing relation use :","possible_answers":"Create table ,Modify table ,Alter table ,Drop table ","correct_answer":"3"} ,{"id":"10","question":"Which of the following statements creates a new table temp instructor that has the same schema as an instructor.","possible_answers":"create table temp_instructor; ,create table temp_instructor like instructor; ,create Table as temp_instructor; ,create table like temp_instructor; ","correct_answer":"2"} ,{"id":"11","question":"Foreign key is the one in which the ________ of one relation is referenced in another relation.","possible_answers":"Foreign key ,Primary key ,References ,Check constraint ","correct_answer":"2"} ,{"id":"12","question":"SELECT emp_name\nFROM department\nWHERE dept_name LIKE ’ _____ Computer Science’;\nWhich one of the following has to be added into the blank to select the dept_name which has Computer Science as its ending string?\n","possible_answers":"% ,_ ,|| ,$ ","correct_answer":"1"} ,{"id":"13","question":"The term attribute refers to a ___________ of a table.","possible_answers":"Record ,Column ,Tuple ,Key ","correct_answer":"2"} ,{"id":"14","question":"For each attribute of a relation, there is a set of permitted values, called the ________ of that attribute.","possible_answers":"Domain ,Relation ,Set ,Schema ","correct_answer":"1"} ,{"id":"15","question":"Dates must be specified in the format","possible_answers":"mm/dd/yy ,yyyy/mm/dd ,dd/mm/yy ,yy/dd/mm ","correct_answer":"2"} ,{"id":"16","question":"CREATE DOMAIN YearlySalary NUMERIC(8,2)\nCONSTRAINT salary VALUE test __________;\nIn order to ensure that an instructor’s salary domain allows only values greater than a specified value use:\n","possible_answers":"Value>=30000.00 ,Not null; ,Check(value >= 29000.00); ,Check(value) ","correct_answer":"3"} ,{"id":"17","question":"Course(course_id,sec_id,semester)\nHere the course_id,sec_id and semester are __________ and course is a _________\n","possible_answers":"Relations/Attribute ,Attributes/Relation ,Tuple/Relation ,Tuple/Attributes ","correct_answer":"2"} ,{"id":"18","question":"SELECT name\nFROM instructor\nWHERE salary IS NOT NULL;\nSelects\n","possible_answers":"Tuples with null value ,Tuples with no null values ,Tuples with any salary ,All of the mentioned ","correct_answer":"2"} ,{"id":"19","question":"In an employee table to include the attributes whose value always have some value which of the following constraint must be used?","possible_answers":"Null ,Not null ,Unique ,Distinct ","correct_answer":"2"} ,{"id":"20","question":"Which one of the following is used to define the structure of the relation, deleting relations and relating schemas?","possible_answers":"DML(Data Manipulation Langauge) ,DDL(Data Definition Langauge) ,Query ,Relational Schema ","correct_answer":"2"} ,{"id":"21","question":"The primary key must be","possible_answers":"Unique ,Not null ,Both Unique and Not null ,Either Unique or Not null ","correct_answer":"3"} ,{"id":"22","question":"How many tables may be included with a join?","possible_answers":"One ,Two ,Three ,All of the mentioned ","correct_answer":"4"} ,{"id":"23","question":"Which one of the following attribute can be taken as a primary key?","possible_answers":"Name ,Street ,Id ,Department ","correct_answer":"3"} ,{"id":"24","question":"The default extension for an Oracle SQL*Plus file is:","possible_answers":".txt ,.pls ,.ora ,.sql ","correct_answer":"4"} ,{"id":"25","question":"Which of the following deletes all tuples in the instructor relation for those instructors associated with a department located in the Watson building which is in department relation.","possible_answers":"DELETE FROM instructor\n   WHERE dept_name IN ''Watson'';\n ,DELETE FROM department \n   WHERE building=''Watson'';\n ,DELETE FROM instructor\nWHERE dept_name IN (SELECT dept name\nFROM department\nWHERE building = ’Watson’);\n ,None of the mentioned ","correct_answer":"3"} ]}'
